$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("C54:H69").NumberFormat = "0.00"
Write-Output "done"
